$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" (changed) date column C was bumped from 2023-09-15 (45184)
# to 2023-09-16 (45185) for rows 2 through 12.
$newDate = [DateTime]::FromOADate(45185)

for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
